$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (Forecasted Load), C (Solar), D (Wind) for rows 2-25
# Row, B, C, D
$data = @(
    @(2, 48061, 0, 37962),
    @(3, 46536, 0, 36335),
    @(4, 45639, 0, 35011),
    @(5, 46060, 0, 33667),
    @(6, 47295, 0, 32154.5),
    @(7, 50633, 0, 30641),
    @(8, 57774, 1, 29213),
    @(9, 62858, 2071, 27597),
    @(10, 65338, 8050, 25805.5),
    @(11, 66158, 15549, 24000),
    @(12, 67056, 21685, 22569.5),
    @(13, 67572, 23653, 21864),
    @(14, 67151, 22393, 21418.5),
    @(15, 65494, 18592, 20716),
    @(16, 63993, 12614, 19977.5),
    @(17, 62602, 6051, 20148),
    @(18, 62301, 1397, 21592),
    @(19, 65596, 5, 23772),
    @(20, 65385, 0, 25521.5),
    @(21, 64085, 0, 26407),
    @(22, 60694, 0, 26297),
    @(23, 57200, 0, 25465.5),
    @(24, 53743, 0, 24760),
    @(25, 49913, 0, 24705)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
}
